$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Beams")

# Termino beam summary de corte (Vz column E): set values to 0, replacing formulas
$ws.Range("E7:E14").Value = 0

# Update selection to I9
$ws.Range("I9").Select()
